$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: paragraph "{m:'zone1'.myTemplate()}" -> insert a brand-new
# run carrying a single space between the opening-quote run and the
# "zone1" run, giving "{m:' zone1'.myTemplate()}".
# ---------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("{m:'zone1'.myTemplate()}")
if ($found1) {
    $target1 = $d.Range($rng1.Start, $rng1.End)
    # NOTE: the "'" and ".myTemplate()}" runs originally carry
    # <w:rPr><w:lang w:val="en-US"/></w:rPr> - that formatting must be
    # reproduced explicitly, otherwise InsertXML silently drops it from
    # every run it rewrites (only the brand-new space run has no rPr,
    # matching the diff).
    $xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>{m:</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>''</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>zone1</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>''</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">.myTemplate()}</w:t></w:r></w:p></w:body></w:document>
</pkg:xmlData></pkg:part></pkg:package>'
    $target1.InsertXML($xml1)
}

# ---------------------------------------------------------------------
# Change 2: paragraph "{m:userdoc str}" -> split the run holding "str}"
# into "str" and "}" so that the pre-existing _GoBack bookmark ends up
# sitting between the two runs instead of after both of them.
#
# The bookmark has to be removed before the surrounding text is
# rewritten (otherwise InsertXML merges runs back together and leaves
# duplicated bookmark markers behind) and then re-created at the right
# spot once the runs are in their final shape.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$rng2 = $d.Content
$found2 = $rng2.Find.Execute("str}")
if ($found2) {
    $target2 = $d.Range($rng2.Start, $rng2.End)
    $xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>str</w:t></w:r><w:r><w:t xml:space="preserve">}</w:t></w:r></w:p></w:body></w:document>
</pkg:xmlData></pkg:part></pkg:package>'
    $target2.InsertXML($xml2)
}

$rng3 = $d.Content
$found3 = $rng3.Find.Execute("{m:userdoc str")
if ($found3) {
    $bmPoint = $d.Range($rng3.End, $rng3.End)
    $d.Bookmarks.Add("_GoBack", $bmPoint)
}
